# "Generate Report for Archive" - regenerated localization-status report.
#
# The report generator re-sorted the three rows that previously covered
# c646216a / ec7f7bc4 / 4ea2a910 so that 4ea2a910 now sorts first: the data
# for each GUID stays intact, but it now lands on a different row:
#   old row 3 (c646216a) -> new row 4
#   old row 4 (ec7f7bc4) -> new row 5
#   old row 5 (4ea2a910) -> new row 3
# This happens identically on the "Overview", "zh-cn" and "de-de" sheets.
# Rows 2, 6 and 7 (77153fb4, 69b6bb8f, abee8654) are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (file name), B (path+name, hyperlinked),
# E/F (status), G (datetime). C/D are identical across rows 3-5 already.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
$wsOverview.Range("B3").Value = "e2e\4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 02:37:59"

$wsOverview.Range("A4").Value = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
$wsOverview.Range("B4").Value = "e2e\c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-08-16 02:37:28"

$wsOverview.Range("A5").Value = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
$wsOverview.Range("B5").Value = "e2e\ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("G5").Value = "2016-08-16 02:37:28"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
    } elseif ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
    } elseif ($addr -eq '$B$5') {
        $hl.TextToDisplay = "e2e\ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (source file, hyperlinked), C (status),
# G (latest handoff file), H (latest handoff datetime).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.47fef88e90054f4f3a31cd9b89f8dd1cebb2be51.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 02:37:53"

$wsZhCn.Range("A4").Value = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.05c25e210db9d40c6b5f749af062eee66a0eaaeb.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-16 02:37:23"

$wsZhCn.Range("A5").Value = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("G5").Value = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.27c681e9b1319bc5bf27fd8a3e15db846eee2634.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-16 02:37:23"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
    } elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
    } elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn".
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.47fef88e90054f4f3a31cd9b89f8dd1cebb2be51.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 02:37:59"

$wsDeDe.Range("A4").Value = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.05c25e210db9d40c6b5f749af062eee66a0eaaeb.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-16 02:37:28"

$wsDeDe.Range("A5").Value = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("G5").Value = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.27c681e9b1319bc5bf27fd8a3e15db846eee2634.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-16 02:37:28"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "4ea2a910-4ddf-4b0d-a217-b1f3aadf60ef.md"
    } elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "c646216a-d336-4d07-a7a5-d1bf5c66dd15.md"
    } elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "ec7f7bc4-68c3-409e-9d60-cf9ee38a1fd6.md"
    }
}
